# Daily attendance processing - 2025-11-24 14:26:32
# Reverse the order of the comma-separated "Recorded By" entries in column G
# for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"
        $n = $parts.Count
        if ($n -gt 1) {
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
